# This script updates the "想去人数" (interest count, column F) values on
# three worksheets of the 杭州-漫展信息 workbook to match a refreshed data
# pull (per commit "Update gh-pages to output generated at 456a3b4").
#
# Sheet names (order in workbook): 展览, 演出, 本地生活, 全部类型

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 9053
$ws1.Range("F4").Value  = 6627
$ws1.Range("F6").Value  = 2148
$ws1.Range("F14").Value = 85
$ws1.Range("F15").Value = 26
$ws1.Range("F16").Value = 8956
$ws1.Range("F28").Value = 1043
$ws1.Range("F31").Value = 560
$ws1.Range("F32").Value = 34
$ws1.Range("F33").Value = 41
$ws1.Range("F34").Value = 550
$ws1.Range("F35").Value = 2365
$ws1.Range("F37").Value = 551
$ws1.Range("F42").Value = 183
$ws1.Range("F46").Value = 85
$ws1.Range("F47").Value = 18
$ws1.Range("F49").Value = 9

# ---------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F14").Value = 19

# ---------------------------------------------------------------
# Sheet "全部类型" (All types)
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 9053
$ws4.Range("F6").Value  = 6627
$ws4.Range("F8").Value  = 2148
$ws4.Range("F17").Value = 85
$ws4.Range("F18").Value = 8956
$ws4.Range("F27").Value = 1043
$ws4.Range("F31").Value = 560
$ws4.Range("F32").Value = 34
$ws4.Range("F33").Value = 41
$ws4.Range("F34").Value = 550
$ws4.Range("F35").Value = 2365
$ws4.Range("F37").Value = 19
$ws4.Range("F39").Value = 551
$ws4.Range("F41").Value = 183
$ws4.Range("F43").Value = 85
